# Update "想去人数" (number of people interested) values on the
# "展览" and "全部类型" worksheets, as published by the gh-pages scraper run.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 1900
    $ws.Range("F3").Value = 358
    $ws.Range("F5").Value = 1225
    $ws.Range("F7").Value = 6008
}
